# Switch licence from BY-NC to BY-SA
#
# The document credits the work under a Creative Commons BY-NC licence;
# this script updates both the human-readable licence name in the body
# text and the hyperlink (address + display text) that points at the
# creativecommons.org licence page, so everything consistently reads
# "BY-SA" instead of "BY-NC".

$d = $word.ActiveDocument

# 1) Update the visible licence text: "CC BY-NC 4.0" -> "CC BY-SA 4.0"
#    (wildcards off, match case on, whole match only — there is exactly
#    one occurrence, in the licensing paragraph near the top of the
#    references list).
$d.Content.Find.Execute("CC BY-NC 4.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CC BY-SA 4.0", 2) | Out-Null

# 2) Update the Creative Commons hyperlink itself so its target and its
#    displayed URL text both point at the by-sa licence instead of by-nc.
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.Address -like "*creativecommons.org/licenses/by-nc/4.0*") {
        $h.Address = "https://creativecommons.org/licenses/by-sa/4.0"
        $h.TextToDisplay = "https://creativecommons.org/licenses/by-sa/4.0"
        break
    }
}
